# Converts an "RRGGBB" hex string into the packed OLE/VBA BGR color
# integer (0x00BBGGRR) that the PowerPoint object model's RGB
# properties expect.
function HexToOle($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Swap the deck's applied theme colours from the "Integral" /
#        "Red Violet" scheme to the default "Office" scheme. ---
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Item($i).RGB = HexToOle $officeColors[$i - 1]
}

# --- 2. Re-style the table on slide 5 with the new table style GUID. ---
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shape = $tableSlide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{0C3319E0-3D4A-4C3A-9FFE-D13112985D67}")
    }
}
